# "Added Import Budget Model, DataSheet , NameMapping"
# Appends a new "ImportBudgetModel" test-case row (row 57) to the
# GlobalTestCase sheet, following the same layout as the existing rows
# (UnitName/TestCases/Description repeat the name, the region columns
# default to "No", and the last column is "Yes").

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("GlobalTestCase")

$ws.Range("A57").Value = "ImportBudgetModel"
$ws.Range("B57").Value = "ImportBudgetModel"
$ws.Range("C57").Value = "ImportBudgetModel"
$ws.Range("D57").Value = "No"
$ws.Range("E57").Value = "No"
$ws.Range("F57").Value = "No"
$ws.Range("G57").Value = "No"
$ws.Range("H57").Value = "No"
$ws.Range("I57").Value = "Yes"

# Match the text number format used by the rest of the table (style index 2).
$ws.Range("A57:I57").NumberFormat = "@"

# Mirror the selection left behind in the saved workbook.
$ws.Range("A57:I57").Select()
